$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" (total) summary sheet -----------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Push the existing 2022-Q2 summary row down to row 3 (copy with formatting).
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3").Value = 1

# Write the new 2022-Q4 summary figures into row 2.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = 0.01

# --- 2. Insert a new "2022-Q4" detail sheet, right before "2022-Q2" --------
$oldQ2Sheet = $wb.Worksheets.Item("2022-Q2")
$q4Sheet = $wb.Worksheets.Add($oldQ2Sheet)
$q4Sheet.Name = "2022-Q4"

# Match the page margins used on the 总计 sheet.
$q4Sheet.PageSetup.LeftMargin = $totalSheet.PageSetup.LeftMargin
$q4Sheet.PageSetup.RightMargin = $totalSheet.PageSetup.RightMargin
$q4Sheet.PageSetup.TopMargin = $totalSheet.PageSetup.TopMargin
$q4Sheet.PageSetup.BottomMargin = $totalSheet.PageSetup.BottomMargin
$q4Sheet.PageSetup.HeaderMargin = $totalSheet.PageSetup.HeaderMargin
$q4Sheet.PageSetup.FooterMargin = $totalSheet.PageSetup.FooterMargin

# Copy header formatting (bold/border style) from the 总计 sheet header cell.
$totalSheet.Range("B1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A3").PasteSpecial(-4122)

# Header row.
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Row 2 data (fund A).
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2:G2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "011886"
$q4Sheet.Range("C2").Value = "弘毅远方高端制造混合型发起式证券投资基金A"
$q4Sheet.Range("D2").Value = "0.23"
$q4Sheet.Range("E2").Value = "88.95"
$q4Sheet.Range("F2").Value = "2.96"
$q4Sheet.Range("G2").Value = "0.0068"
$q4Sheet.Range("H2").Value = 10

# Row 3 data (fund C).
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3:G3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "011887"
$q4Sheet.Range("C3").Value = "弘毅远方高端制造混合型发起式证券投资基金C"
$q4Sheet.Range("D3").Value = "0.12"
$q4Sheet.Range("E3").Value = "88.95"
$q4Sheet.Range("F3").Value = "2.96"
$q4Sheet.Range("G3").Value = "0.0036"
$q4Sheet.Range("H3").Value = 10

# Keep 2022-Q2 as the active/selected sheet, matching original selection.
$oldQ2Sheet.Activate()
